$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1088.2
$ws.Range("I18").Value = 1111.25
$ws.Range("J18").Value = 996
$ws.Range("K18").Value = 1111.25
$ws.Range("L18").Value = 996
$ws.Range("M18").Value = -827.25
$ws.Range("N18").Value = -1564
$ws.Range("H40").Value = 4600.1665
$ws.Range("I40").Value = 4100
$ws.Range("J40").Value = 5600.5
$ws.Range("K40").Value = 4100
$ws.Range("L40").Value = 5600.5
$ws.Range("M40").Value = -3925
$ws.Range("N40").Value = -5950.5
$ws.Range("H41").Value = 3749.6667
$ws.Range("I41").Value = 3500
$ws.Range("J41").Value = 3874.5
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 3874.5
$ws.Range("M41").Value = -3060
$ws.Range("N41").Value = -4754.5
$ws.Range("H42").Value = 130
$ws.Range("I42").Value = 130
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 390
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -160
$ws.Range("N42").Value = ""
$ws.Range("H96").Value = 206.53847
$ws.Range("J96").Value = 207.66667
$ws.Range("L96").Value = 623.00001
$ws.Range("N96").Value = -3369.00001
$ws.Range("H113").Value = 5224
$ws.Range("I113").Value = 7600
$ws.Range("J113").Value = 3739
$ws.Range("K113").Value = 7600
$ws.Range("L113").Value = 3739
$ws.Range("M113").Value = -4346
$ws.Range("N113").Value = -10247
$ws.Range("H137").Value = 3574.6875
$ws.Range("I137").Value = 2248
$ws.Range("J137").Value = 3880.8462
$ws.Range("K137").Value = 6744
$ws.Range("L137").Value = 11642.5386
$ws.Range("M137").Value = -4194
$ws.Range("N137").Value = -16742.5386

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 6042930
$ws.Range("I43").Value = 7528662
$ws.Range("K43").Value = 7528662
$ws.Range("M43").Value = -7528349
$ws.Range("H61").Value = 4663.3335
$ws.Range("I61").Value = 4663.3335
$ws.Range("K61").Value = 4663.3335
$ws.Range("M61").Value = -4451.3335
$ws.Range("H112").Value = 11500
$ws.Range("J112").Value = 11500
$ws.Range("L112").Value = 11500
$ws.Range("N112").Value = -14454
$ws.Range("H114").Value = 49999
$ws.Range("J114").Value = 49999
$ws.Range("L114").Value = 49999
$ws.Range("N114").Value = -58677
$ws.Range("H122").Value = 2386.5833
$ws.Range("J122").Value = 3231.3333
$ws.Range("L122").Value = 9693.999899999999
$ws.Range("N122").Value = -14593.9999
$ws.Range("H136").Value = 4663.3335
$ws.Range("I136").Value = 4663.3335
$ws.Range("K136").Value = 13990.0005
$ws.Range("M136").Value = -11440.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 15083
$ws.Range("I54").Value = 15083
$ws.Range("K54").Value = 15083
$ws.Range("M54").Value = -14599
$ws.Range("H112").Value = 46049.8
$ws.Range("J112").Value = 46049.8
$ws.Range("L112").Value = 46049.8
$ws.Range("N112").Value = -49003.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7676.6924
$ws.Range("I31").Value = 4205.7144
$ws.Range("J31").Value = 8955.474
$ws.Range("K31").Value = 4205.7144
$ws.Range("L31").Value = 8955.474
$ws.Range("M31").Value = -3910.7144
$ws.Range("N31").Value = -9545.474
$ws.Range("H34").Value = 7676.6924
$ws.Range("I34").Value = 4205.7144
$ws.Range("J34").Value = 8955.474
$ws.Range("K34").Value = 4205.7144
$ws.Range("L34").Value = 8955.474
$ws.Range("M34").Value = -4003.7144
$ws.Range("N34").Value = -9359.474
$ws.Range("H58").Value = 4691.6665
$ws.Range("I58").Value = 1252.2
$ws.Range("K58").Value = 1252.2
$ws.Range("M58").Value = -1049.2
$ws.Range("H105").Value = 2144
$ws.Range("I105").Value = 2144
$ws.Range("K105").Value = 2144
$ws.Range("M105").Value = -397
$ws.Range("H136").Value = 4691.6665
$ws.Range("I136").Value = 1252.2
$ws.Range("K136").Value = 3756.6
$ws.Range("M136").Value = -1206.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2100
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 2375
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 7125
$ws.Range("M5").Value = -2888
$ws.Range("N5").Value = -7349
$ws.Range("H49").Value = 2500
$ws.Range("I49").Value = 2500
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 7500
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -7344
$ws.Range("N49").Value = ""
$ws.Range("H135").Value = 2100
$ws.Range("I135").Value = 1000
$ws.Range("J135").Value = 2375
$ws.Range("K135").Value = 9000
$ws.Range("L135").Value = 21375
$ws.Range("M135").Value = -6465
$ws.Range("N135").Value = -26445

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1104.5
$ws.Range("I102").Value = 1191.762
$ws.Range("J102").Value = 493.66666
$ws.Range("K102").Value = 1191.762
$ws.Range("L102").Value = 493.66666
$ws.Range("M102").Value = 430.2380000000001
$ws.Range("N102").Value = -3737.66666
$ws.Range("H132").Value = 38199.332
$ws.Range("I132").Value = 47787.74
$ws.Range("K132").Value = 143363.22
$ws.Range("M132").Value = -140833.22

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5895.385
$ws.Range("I7").Value = 4136.5713
$ws.Range("J7").Value = 7947.3335
$ws.Range("K7").Value = 4136.5713
$ws.Range("L7").Value = 7947.3335
$ws.Range("M7").Value = -4024.5713
$ws.Range("N7").Value = -8171.3335
$ws.Range("H22").Value = 1051.9333
$ws.Range("I22").Value = 995.4
$ws.Range("J22").Value = 1080.2
$ws.Range("K22").Value = 995.4
$ws.Range("L22").Value = 1080.2
$ws.Range("M22").Value = -700.4
$ws.Range("N22").Value = -1670.2
$ws.Range("H27").Value = 1051.9333
$ws.Range("I27").Value = 995.4
$ws.Range("J27").Value = 1080.2
$ws.Range("K27").Value = 995.4
$ws.Range("L27").Value = 1080.2
$ws.Range("M27").Value = -888.4
$ws.Range("N27").Value = -1294.2
$ws.Range("H46").Value = 4188.643
$ws.Range("I46").Value = 4718.143
$ws.Range("J46").Value = 4012.1428
$ws.Range("K46").Value = 4718.143
$ws.Range("L46").Value = 4012.1428
$ws.Range("M46").Value = -4530.143
$ws.Range("N46").Value = -4388.1428
$ws.Range("H53").Value = 10637.5
$ws.Range("I53").Value = 10637.5
$ws.Range("K53").Value = 10637.5
$ws.Range("M53").Value = -10119.5
$ws.Range("H64").Value = 22499.75
$ws.Range("J64").Value = 22499.75
$ws.Range("L64").Value = 22499.75
$ws.Range("N64").Value = -22949.75
$ws.Range("H67").Value = 22499.75
$ws.Range("J67").Value = 22499.75
$ws.Range("L67").Value = 22499.75
$ws.Range("N67").Value = -24059.75
$ws.Range("H110").Value = 18875
$ws.Range("J110").Value = 18875
$ws.Range("L110").Value = 18875
$ws.Range("N110").Value = -27055
$ws.Range("H126").Value = 5895.385
$ws.Range("I126").Value = 4136.5713
$ws.Range("J126").Value = 7947.3335
$ws.Range("K126").Value = 12409.7139
$ws.Range("L126").Value = 23842.0005
$ws.Range("M126").Value = -9939.713899999999
$ws.Range("N126").Value = -28782.0005
